# Israel Premier League workbook update
# Swaps betting-odds rows whose underlying match records were re-ordered
# during a re-scrape, and refreshes a handful of odds cells for upcoming
# fixtures (rows 286-291) whose prices moved before kickoff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Simple pairwise swaps: the full record (columns B..AC) of row A and
#    row B trade places; column A (the running index) stays put.
# ---------------------------------------------------------------------
$swapPairs = @(
    @(19, 20),
    @(41, 42),
    @(48, 49),
    @(83, 84),
    @(104, 105),
    @(110, 111),
    @(116, 118),
    @(120, 121),
    @(127, 128),
    @(129, 131),
    @(140, 141),
    @(148, 149),
    @(161, 162),
    @(216, 217),
    @(239, 240)
)

foreach ($pair in $swapPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# ---------------------------------------------------------------------
# 2) Three-way rotations: row0 <- row2, row1 <- row0, row2 <- row1
# ---------------------------------------------------------------------
$rotations = @(
    @(96, 97, 98),
    @(188, 189, 190)
)

foreach ($rot in $rotations) {
    $row0 = $rot[0]
    $row1 = $rot[1]
    $row2 = $rot[2]

    $range0 = $ws.Range("B$row0`:AC$row0")
    $range1 = $ws.Range("B$row1`:AC$row1")
    $range2 = $ws.Range("B$row2`:AC$row2")

    $values0 = $range0.Value()
    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range0.Value = $values2
    $range1.Value = $values0
    $range2.Value = $values1
}

# ---------------------------------------------------------------------
# 3) Direct odds refreshes on upcoming (not-yet-played) fixtures - only
#    specific cells change, no row is swapped.
# ---------------------------------------------------------------------
$ws.Range("N286").Value = 2.8
$ws.Range("P286").Value = 2.4
$ws.Range("Q286").Value = 0
$ws.Range("R286").Value = 2.1
$ws.Range("S286").Value = 1.775
$ws.Range("T286").Value = 2.5
$ws.Range("U286").Value = 2.025
$ws.Range("V286").Value = 1.825

$ws.Range("U287").Value = 1.875
$ws.Range("V287").Value = 1.975

$ws.Range("R289").Value = 2.025
$ws.Range("S289").Value = 1.825
$ws.Range("U289").Value = 1.875
$ws.Range("V289").Value = 1.975

$ws.Range("O290").Value = 3.5
$ws.Range("R290").Value = 1.825
$ws.Range("S290").Value = 2.025
$ws.Range("T290").Value = 2.5
$ws.Range("U290").Value = 1.925
$ws.Range("V290").Value = 1.925

$ws.Range("R291").Value = 1.9
$ws.Range("S291").Value = 1.95
